$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.580.92'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '2.519.59'
$ws.Range("E3").Value = '  -2.99%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '''309.43'
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("D6").Value = '''100.48'
$ws.Range("E6").Value = '  +2.36%  '
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("D10").Value = '''35.92'
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").Value = '''0.0801'
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").Value = '''7.28'
$ws.Range("E12").Value = '  -3.52%  '
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '2.907.40'
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").Value = '2.488.46'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '''0.806'
$ws.Range("E17").Value = '  -5.14%  '
$ws.Range("D18").Value = '42.549.95'
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("D19").Value = '''6.71'
$ws.Range("E19").Value = '  -2.50%  '
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("D21").Value = '''12.10'
$ws.Range("E21").Value = '  -5.33%  '
$ws.Range("D22").Value = '''69.30'
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("D23").Value = '''243.41'
$ws.Range("E23").Value = '  -4.48%  '
$ws.Range("D24").Value = '''2.88'
$ws.Range("E24").Value = '  -3.58%  '
$ws.Range("D25").Value = '''2.03'
$ws.Range("E25").Value = '  -2.88%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").Value = '''25.91'
$ws.Range("E27").Value = '  -5.14%  '
$ws.Range("E28").Value = '  -3.77%  '
$ws.Range("D29").Value = '''10.12'
$ws.Range("D30").Value = '''38.84'
$ws.Range("E30").Value = '  -5.38%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '''5.75'
$ws.Range("E31").Value = '  -2.23%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '''155.48'
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("E33").Value = '  +10.22%  '
$ws.Range("D34").Value = '''0.0785'
$ws.Range("E34").Value = '  -3.20%  '
$ws.Range("E35").Value = '  -3.03%  '
$ws.Range("D36").Value = '''3.19'
$ws.Range("E36").Value = '  -8.23%  '
$ws.Range("E37").Value = '  -6.73%  '
$ws.Range("D38").Value = '''18.21'
$ws.Range("E38").Value = '  -3.42%  '
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("E41").Value = '  +5.63%  '
$ws.Range("D42").Value = '''22.08'
$ws.Range("E42").Value = '  -4.43%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  -2.20%  '
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("D46").Value = '1.982.61'
$ws.Range("E46").Value = '  -1.55%  '
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").Value = '2.762.75'
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("D49").Value = '''80.06'
$ws.Range("E49").Value = '  -4.34%  '
$ws.Range("E50").Value = '  -3.15%  '
$ws.Range("D51").Value = '''72.26'
$ws.Range("E51").Value = '  -3.54%  '
